# Update UG DG diagrams to reflect Event Manager
#
# LogicComponentSequenceDiagram.pptx, slide 1:
#   - The ":Address / BookParser" lifeline label becomes
#     ":EventManager / BookParser" (and shrinks from 16pt to 12pt to
#     keep fitting the lifeline header box).
#   - The "deletePerson(p)" call-out becomes "deleteEvent(e)".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $targetId) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $candidate = $slide.Shapes.Item($i)
        if ($candidate.Id -eq $targetId) { return $candidate }
    }
    return $null
}

# "Rectangle 62" lifeline header: ":Address" / "BookParser" -> ":EventManager" / "BookParser"
$addressBox = Get-ShapeById $s 16
$lifelineText = $addressBox.TextFrame.TextRange
$lifelineText.Text = ":EventManager`rBookParser"
$lifelineText.Font.Size = 12

# "TextBox 77" call-out: "deletePerson(p)" -> "deleteEvent(e)"
$deleteCallout = Get-ShapeById $s 78
$calloutText = $deleteCallout.TextFrame.TextRange
$calloutText.Runs(1, 1).Text = "deleteEvent"
$calloutText.Runs(2, 1).Text = "(e)"
